$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A91").Value = "2025-04-29 15:18:54"
$ws.Range("B91").Value = 245
